# Apply the odds updates described in the commit diff for
# "Jogos_da_Semana_FlashScore_2025-02-19.xlsx"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("N2").Value = 15
$ws.Range("S2").Value = 1.92
$ws.Range("T2").Value = 1.98

# Row 3
$ws.Range("J3").Value = 2.63
$ws.Range("Z3").Value = 1.67

# Row 4
$ws.Range("J4").Value = 2.38
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 9.5
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.8
$ws.Range("Y4").Value = 1.91
$ws.Range("Z4").Value = 1.8
$ws.Range("AB4").Value = 7.5

# Row 5
$ws.Range("G5").Value = 5.5
$ws.Range("I5").Value = 1.57
$ws.Range("Y5").Value = 1.83
$ws.Range("Z5").Value = 1.83
$ws.Range("AK5").Value = 301

# Row 6
$ws.Range("G6").Value = 2.6
$ws.Range("I6").Value = 2.7
$ws.Range("L6").Value = 3.25
$ws.Range("AO6").Value = 26
